$d = $word.ActiveDocument

# --- 1) Paragraph with the "m:if" field: replace the fldChar/instrText
#        run sequence with plain w:t runs wrapped in "{ ... }" braces.
$pIf = $d.Paragraphs(2)
$rIf = $pIf.Range
$xmlIf = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidP="002A1F2A" w:rsidR="002A1F2A" w:rsidRDefault="002A1F2A"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">{m:if </w:t></w:r><w:r><w:t xml:space="preserve">self.name </w:t></w:r><w:r><w:t>&lt;&gt;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>'</w:t></w:r><w:r><w:t>anydsl</w:t></w:r><w:r><w:t>'}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rIf.InsertXML($xmlIf) | Out-Null

# --- 2) Paragraph with the "m:elseif" field: replace the fldChar/instrText
#        run sequence with plain w:t runs wrapped in "{ ... }" braces,
#        keeping the trailing (already plain-text) runs unchanged.
$pElseif = $d.Paragraphs(4)
$rElseif = $pElseif.Range
$xmlElseif = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidP="00430772" w:rsidR="00430772" w:rsidRDefault="00430772"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr><w:r><w:t>{m:</w:t></w:r><w:r><w:t>elseif self.name = 'anydsl'</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>m:elseif, m:else or m:endif expected here while parsing m:elseif self.name = 'anydsl'</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rElseif.InsertXML($xmlElseif) | Out-Null

# --- 3) Last paragraph: drop the trailing run of 4 spaces that followed
#        the "Unexpected tag EOF ..." text.
$pLast = $d.Paragraphs($d.Paragraphs.Count)
$endLast = $pLast.Range.End
$rTrailingSpaces = $d.Range($endLast - 5, $endLast - 1)
$rTrailingSpaces.Delete() | Out-Null
